$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sprint 2")
$ws3 = $wb.Worksheets.Item("Sprint 3")

# ---------------------------------------------------------------------------
# Sprint 3 sheet: fill in "Story Points Completed" (F) and contribution (G)
# for the stories that wrapped up, and append a new backlog row for
# "75, Finalize sprint 3 backlog".
# ---------------------------------------------------------------------------
$ws3.Range("F10").Value = 2
$ws3.Range("F11").Value = 2

$ws3.Range("G12").Value = "Vasilis: 100%"

$ws3.Range("F16").Value = 4
$ws3.Range("G16").Value = "Michael, Sakshyam, Vasilis: 100%"

$ws3.Range("F17").Value = 4
$ws3.Range("G17").Value = "Michael, Sakshyam, Vasilis: 100%"

$ws3.Range("F20").Value = 2
$ws3.Range("F21").Value = 2

# New row 22: "75, Finalize sprint 3 backlog"
$ws3.Range("A22").Value = "75, Finalize sprint 3 backlog"
$ws3.Range("B22").Value = "N/A"
$ws3.Range("C22").Value = "N/A"
$ws3.Range("D22").Value = "T"
$ws3.Range("E22").Value = 1
$ws3.Range("F22").Value = 1
$ws3.Range("G22").Value = "Brody: 100%"

# ---------------------------------------------------------------------------
# Sprint 2 sheet: fix the "In progess" typo -> "In progress" in column F
# for rows 20-23 (Story Points Completed column holding a status note
# instead of a number for these in-flight stories).
# ---------------------------------------------------------------------------
$fixedText = "In progress as of end of sprint "
$ws2.Range("F20").Value = $fixedText
$ws2.Range("F21").Value = $fixedText
$ws2.Range("F22").Value = $fixedText
$ws2.Range("F23").Value = $fixedText

# ---------------------------------------------------------------------------
# View state: Sprint 3 becomes the active/selected tab, Sprint 2 keeps a
# lingering selection at F23 (no longer the active tab).
# ---------------------------------------------------------------------------
$ws2.Range("F23").Select()

$ws3.Activate()
$excel.ActiveWindow.Zoom = 131
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
$ws3.Range("H18").Select()
